$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps 2 RGB")

# Ordered list of new rows (105-118), each as: row number, label, optional B:I values.
# Order matters: it must match the order new shared strings were originally created in
# (V7.04.5 .. V7.04.18, rows 105..118 strictly ascending).
$rows = @(
    @{ Row = 105; Label = "V7.04.5" },
    @{ Row = 106; Label = "V7.04.6" },
    @{ Row = 107; Label = "V7.04.7" },
    @{ Row = 108; Label = "V7.04.8" },
    @{ Row = 109; Label = "V7.04.9" },
    @{ Row = 110; Label = "V7.04.10"; Values = @(7.3808999999999996, 0.65229999999999999, 9.6631999999999998, 0.75600000000000001, 9.4557000000000002, 0.71309999999999996, 8.0066000000000006, 0.72050000000000003) },
    @{ Row = 111; Label = "V7.04.11"; Values = @(7.9814999999999996, 0.66700000000000004, 8.5655999999999999, 0.75060000000000004, 9.9532000000000007, 0.73770000000000002, 7.6144999999999996, 0.71730000000000005) },
    @{ Row = 112; Label = "V7.04.12"; Values = @(9.7713000000000001, 0.63939999999999997, 11.5276, 0.77129999999999999, 9.5265000000000004, 0.71609999999999996, 8.4581, 0.73160000000000003) },
    @{ Row = 113; Label = "V7.04.13"; Values = @(5.9819000000000004, 0.66390000000000005, 8.2500999999999998, 0.74739999999999995, 9.6128, 0.71389999999999998, 7.6289999999999996, 0.71519999999999995) },
    @{ Row = 114; Label = "V7.04.14"; Values = @(9.1943000000000001, 0.65210000000000001, 9.8124000000000002, 0.75109999999999999, 8.2829999999999995, 0.71699999999999997, 8.2120999999999995, 0.72109999999999996) },
    @{ Row = 115; Label = "V7.04.15" },
    @{ Row = 116; Label = "V7.04.16" },
    @{ Row = 117; Label = "V7.04.17" },
    @{ Row = 118; Label = "V7.04.18" }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Label
    if ($entry.ContainsKey("Values")) {
        $vals = $entry.Values
        for ($c = 0; $c -lt $vals.Length; $c++) {
            $ws.Cells.Item($r, $c + 2).Value = $vals[$c]
        }
    }
}

# Match the author's final selection/scroll position (frozen pane scrolled so
# row 107 is the first visible row below the freeze, cursor on H119).
[void]$ws.Activate()
[void]$ws.Range("H119").Select()
